$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# League data refresh (09-03-2024 13:07):
#  - Match rows 9 & 10 (ids 6814330 / 6814328) had their fixture/odds data
#    swapped between the two rows.
#  - Rows 118-121 are refreshed with updated match data (row that used to be
#    row 122 is gone; each row's odds/id/date/teams are updated in place).
#  - Row 122 is removed (the schedule shrank by one fixture).

# Row 9
$ws.Range("AA9").Value = 0.4375
$ws.Range("AB9").Value = -0.5
$ws.Range("AC9").Value = 0.425
$ws.Range("B9").Value = 6814328
$ws.Range("G9").Value = "NK Bravo"
$ws.Range("I9").Value = 1
$ws.Range("J9").Value = "D"
$ws.Range("K9").Value = 2.35
$ws.Range("L9").Value = 3.1
$ws.Range("M9").Value = 2.9
$ws.Range("N9").Value = 2.15
$ws.Range("O9").Value = 3.1
$ws.Range("P9").Value = 3.3
$ws.Range("Q9").Value = -0.25
$ws.Range("R9").Value = 1.925
$ws.Range("S9").Value = 1.875
$ws.Range("T9").Value = 2.25
$ws.Range("U9").Value = 1.95
$ws.Range("V9").Value = 1.85
$ws.Range("W9").Value = -1
$ws.Range("X9").Value = 2.1
# Row 10
$ws.Range("AA10").Value = 0.475
$ws.Range("AB10").Value = -1
$ws.Range("AC10").Value = 1
$ws.Range("B10").Value = 6814330
$ws.Range("G10").Value = "NK Aluminij"
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = "H"
$ws.Range("K10").Value = 1.363
$ws.Range("L10").Value = 4.5
$ws.Range("M10").Value = 7
$ws.Range("N10").Value = 1.4
$ws.Range("O10").Value = 4.5
$ws.Range("P10").Value = 7
$ws.Range("Q10").Value = -1.25
$ws.Range("R10").Value = 1.85
$ws.Range("S10").Value = 1.95
$ws.Range("T10").Value = 2.75
$ws.Range("U10").Value = 1.8
$ws.Range("V10").Value = 2
$ws.Range("W10").Value = 0.3999999999999999
$ws.Range("X10").Value = -1
# Row 118
$ws.Range("B118").Value = 6814423
$ws.Range("E118").Value = 45360.5625
$ws.Range("F118").Value = "NK Celje"
$ws.Range("G118").Value = "NK Maribor"
$ws.Range("K118").Value = 2
$ws.Range("L118").Value = 3.3
$ws.Range("M118").Value = 3.5
$ws.Range("N118").Value = 2.1
$ws.Range("O118").Value = 3.25
$ws.Range("P118").Value = 3.25
$ws.Range("Q118").Value = -0.25
$ws.Range("U118").Value = 1.925
$ws.Range("V118").Value = 1.875
# Row 119
$ws.Range("B119").Value = 6816449
$ws.Range("E119").Value = 45361.375
$ws.Range("F119").Value = "NK Rogaska"
$ws.Range("G119").Value = "NK Radomlje"
$ws.Range("K119").Value = 2.625
$ws.Range("L119").Value = 3.2
$ws.Range("M119").Value = 2.5
$ws.Range("N119").Value = 2.5
$ws.Range("O119").Value = 3.2
$ws.Range("P119").Value = 2.625
$ws.Range("Q119").Value = 0
$ws.Range("R119").Value = 1.85
$ws.Range("S119").Value = 1.95
$ws.Range("U119").Value = 2
$ws.Range("V119").Value = 1.8
# Row 120
$ws.Range("B120").Value = 6814420
$ws.Range("E120").Value = 45361.45833333334
$ws.Range("F120").Value = "NS Mura"
$ws.Range("G120").Value = "Olimpija Ljubljana"
$ws.Range("K120").Value = 5.75
$ws.Range("L120").Value = 4
$ws.Range("M120").Value = 1.5
$ws.Range("N120").Value = 5.25
$ws.Range("O120").Value = 3.8
$ws.Range("P120").Value = 1.55
$ws.Range("Q120").Value = 1
$ws.Range("R120").Value = 1.85
$ws.Range("S120").Value = 1.95
$ws.Range("U120").Value = 1.825
$ws.Range("V120").Value = 1.975
# Row 121
$ws.Range("B121").Value = 6814422
$ws.Range("E121").Value = 45361.67708333334
$ws.Range("F121").Value = "NK Domzale"
$ws.Range("G121").Value = "FC Koper"
$ws.Range("K121").Value = 3
$ws.Range("L121").Value = 3.25
$ws.Range("M121").Value = 2.2
$ws.Range("N121").Value = 3.4
$ws.Range("O121").Value = 3.3
$ws.Range("P121").Value = 2
$ws.Range("Q121").Value = 0.25
$ws.Range("R121").Value = 2
$ws.Range("S121").Value = 1.8
$ws.Range("U121").Value = 1.9
$ws.Range("V121").Value = 1.9

# Remove row 122 entirely (shifts nothing else; last row of data)
$ws.Rows("122:122").Delete()
